$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the project number/name values from column B (rows 1-2) to column C
$ws.Range("C1").Value = $ws.Range("B1").Value2
$ws.Range("C2").Value = $ws.Range("B2").Value2
$ws.Range("B1").ClearContents()
$ws.Range("B2").ClearContents()

# Update formulas in column C (rows 5-33) to reference $C$1 instead of $B$1
$ws.Range("C5").Formula = '=_xlfn.CONCAT($C$1, "-",TEXT(A5, "0000"))'
for ($r = 6; $r -le 33; $r++) {
    $ws.Cells.Item($r, 3).Formula = '=IF(B' + $r + '="","",_xlfn.CONCAT($C$1, "-",TEXT(A' + $r + ', "0000")))'
}

# Adjust column A width slightly
$ws.Columns.Item(1).ColumnWidth = 13.82

# Update the selected cell
$ws.Range("E31").Select()
